# "Fragen.xlsx" upload update:
# The attribute spec for the "Hast du ein Auto?" question (row 2, column E)
# now marks "Marke" and "Baujahr" as mandatory ("pflicht"); "Farbe" stays optional.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E2").Value = "Marke:text:pflicht; Farbe:text; Baujahr:number:pflicht"

# Leave the sheet with E2 selected, matching where the author ended up editing.
$ws.Range("E2").Select()

Write-Output "Updated E2 attribute spec and selected E2."
